$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '41.075.06'
$ws.Cells.Item(2, 5).Value = '  -2.10%  '
$ws.Cells.Item(3, 4).Value = '2.159.30'
$ws.Cells.Item(3, 5).Value = '  -2.52%  '
$ws.Cells.Item(4, 5).Value = '  -0.07%  '
$ws.Cells.Item(5, 4).Value = '235.84'
$ws.Cells.Item(5, 5).Value = '  -1.96%  '
$ws.Cells.Item(6, 5).Value = '  -3.35%  '
$ws.Cells.Item(7, 4).Value = '69.08'
$ws.Cells.Item(7, 5).Value = '  -5.65%  '
$ws.Cells.Item(8, 5).Value = '  +0.00%  '
$ws.Cells.Item(9, 4).Value = '0.565'
$ws.Cells.Item(9, 5).Value = '  -6.98%  '
$ws.Cells.Item(10, 4).Value = '38.64'
$ws.Cells.Item(10, 5).Value = '  -9.01%  '
$ws.Cells.Item(11, 4).Value = '0.0909'
$ws.Cells.Item(11, 5).Value = '  -4.37%  '
$ws.Cells.Item(12, 4).Value = '54.25'
$ws.Cells.Item(12, 5).Value = '  -5.65%  '
$ws.Cells.Item(13, 4).Value = "'0.1000"
$ws.Cells.Item(13, 5).Value = '  -2.44%  '
$ws.Cells.Item(14, 4).Value = '6.64'
$ws.Cells.Item(14, 5).Value = '  -6.21%  '
$ws.Cells.Item(15, 4).Value = '2.475.79'
$ws.Cells.Item(15, 5).Value = '  -2.83%  '
$ws.Cells.Item(16, 4).Value = '14.16'
$ws.Cells.Item(16, 5).Value = '  -0.76%  '
$ws.Cells.Item(17, 4).Value = '2.143.62'
$ws.Cells.Item(17, 5).Value = '  -3.45%  '
$ws.Cells.Item(18, 4).Value = '0.781'
$ws.Cells.Item(18, 5).Value = '  -6.66%  '
$ws.Cells.Item(19, 4).Value = '40.815.22'
$ws.Cells.Item(20, 4).Value = '0.0₃0987'
$ws.Cells.Item(20, 5).Value = '  -7.31%  '
$ws.Cells.Item(21, 4).Value = '69.54'
$ws.Cells.Item(21, 5).Value = '  -4.85%  '
$ws.Cells.Item(22, 4).Value = '5.76'
$ws.Cells.Item(22, 5).Value = '  -6.44%  '
$ws.Cells.Item(23, 4).Value = '224.14'
$ws.Cells.Item(23, 5).Value = '  -2.24%  '
$ws.Cells.Item(24, 4).Value = '9.26'
$ws.Cells.Item(24, 5).Value = '  -14.34%  '
$ws.Cells.Item(25, 5).Value = '  +0.03%  '
$ws.Cells.Item(26, 5).Value = '  -10.02%  '
$ws.Cells.Item(27, 4).Value = '10.59'
$ws.Cells.Item(27, 5).Value = '  -9.57%  '
$ws.Cells.Item(28, 4).Value = '3.43'
$ws.Cells.Item(28, 5).Value = '  -5.13%  '
$ws.Cells.Item(29, 5).Value = '  -4.06%  '
$ws.Cells.Item(30, 5).Value = '  -1.28%  '
$ws.Cells.Item(31, 4).Value = '167.86'
$ws.Cells.Item(31, 5).Value = '  +0.61%  '
$ws.Cells.Item(32, 4).Value = '19.74'
$ws.Cells.Item(32, 5).Value = '  -3.44%  '
$ws.Cells.Item(33, 4).Value = '29.98'
$ws.Cells.Item(33, 5).Value = '  +2.22%  '
$ws.Cells.Item(34, 4).Value = '0.0748'
$ws.Cells.Item(34, 5).Value = '  -5.76%  '
$ws.Cells.Item(35, 4).Value = '5.07'
$ws.Cells.Item(35, 5).Value = '  -9.41%  '
$ws.Cells.Item(36, 5).Value = '  -4.23%  '
$ws.Cells.Item(37, 4).Value = '0.101'
$ws.Cells.Item(37, 5).Value = '  -8.42%  '
$ws.Cells.Item(38, 4).Value = '4.05'
$ws.Cells.Item(38, 5).Value = '  -5.06%  '
$ws.Cells.Item(39, 4).Value = '0.0278'
$ws.Cells.Item(39, 5).Value = '  -6.88%  '
$ws.Cells.Item(40, 4).Value = '2.04'
$ws.Cells.Item(40, 5).Value = '  -3.64%  '
$ws.Cells.Item(41, 4).Value = '11.48'
$ws.Cells.Item(41, 5).Value = '  -16.54%  '
$ws.Cells.Item(42, 4).Value = '5.28'
$ws.Cells.Item(42, 5).Value = '  -5.98%  '
$ws.Cells.Item(43, 4).Value = '57.84'
$ws.Cells.Item(43, 5).Value = '  -12.01%  '
$ws.Cells.Item(44, 4).Value = '0.186'
$ws.Cells.Item(44, 5).Value = '  -5.61%  '
$ws.Cells.Item(45, 4).Value = '8.19'
$ws.Cells.Item(45, 5).Value = '  -5.83%  '
$ws.Cells.Item(46, 4).Value = '0.0954'
$ws.Cells.Item(46, 5).Value = '  -4.68%  '
$ws.Cells.Item(47, 4).Value = '96.33'
$ws.Cells.Item(47, 5).Value = '  -7.46%  '
$ws.Cells.Item(48, 4).Value = '1.07'
$ws.Cells.Item(48, 5).Value = '  -4.58%  '
$ws.Cells.Item(49, 5).Value = '  -5.10%  '
$ws.Cells.Item(50, 2).Value = 'NEARProtocol'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(50, 4).Value = '2.15'
$ws.Cells.Item(50, 5).Value = '  -9.51%  '
$ws.Cells.Item(51, 2).Value = 'HuobiToken'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(51, 4).Value = '2.61'
$ws.Cells.Item(51, 5).Value = '  -3.29%  '
